# centrality.xlsx recompute/refresh:
#   - the "name" column (G) for a few rows was re-paired with a different
#     record (the underlying per-character stats moved to a different row
#     while the simple row index in column A stayed positional), and
#   - a number of eigenvector_centrality / betweenness_centrality values
#     were refreshed with slightly different floating point results.
# Columns: A=index, B=betweenness_centrality, C=degree, D=degree_centrality,
#          E=eigenvector_centrality, F=in_degree, G=name, H=out_degree

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 3: tiny eigenvector_centrality refresh
$ws.Cells.Item(3, 5).Value = 0.05445675274586533

# rows 4 & 18: "Zett Jukassa" (row 4) and "Ask Aak" (row 18) swap places
$ws.Cells.Item(4, 2).Value = 0.0003405188414178525
$ws.Cells.Item(4, 5).Value = 0.005073025570724388
$ws.Cells.Item(4, 6).Value = 2
$ws.Cells.Item(4, 7).Value = "Ask Aak"
$ws.Cells.Item(4, 8).Value = 5

$ws.Cells.Item(18, 2).Value = 0.006255618819298771
$ws.Cells.Item(18, 5).Value = 0.009214527369336327
$ws.Cells.Item(18, 6).Value = 3
$ws.Cells.Item(18, 7).Value = "Zett Jukassa"
$ws.Cells.Item(18, 8).Value = 4

# scattered tiny eigenvector_centrality / betweenness_centrality refreshes
$ws.Cells.Item(25, 5).Value = 0.05460153584238549
$ws.Cells.Item(27, 5).Value = 0.1698786606999853
$ws.Cells.Item(41, 2).Value = 0.066474623409301
$ws.Cells.Item(44, 2).Value = 0.007349445002843251
$ws.Cells.Item(45, 2).Value = 0.03694780541905423
$ws.Cells.Item(48, 5).Value = 0.01067139663817989
$ws.Cells.Item(49, 5).Value = 0.01067139663817989
$ws.Cells.Item(55, 5).Value = 0.08291215631090615
$ws.Cells.Item(61, 2).Value = 0.009797724539961626
$ws.Cells.Item(71, 2).Value = 0.008744531328441585

# rows 87 & 89: "Teta (Duke)" (row 87) and "Wat Tambor" (row 89) swap places
$ws.Cells.Item(87, 2).Value = 0.006641473870964542
$ws.Cells.Item(87, 3).Value = 19
$ws.Cells.Item(87, 4).Value = 0.1338028169014085
$ws.Cells.Item(87, 5).Value = 0.06912731264110175
$ws.Cells.Item(87, 6).Value = 8
$ws.Cells.Item(87, 7).Value = "Wat Tambor"
$ws.Cells.Item(87, 8).Value = 11

$ws.Cells.Item(89, 2).Value = 0
$ws.Cells.Item(89, 3).Value = 8
$ws.Cells.Item(89, 4).Value = 0.05633802816901409
$ws.Cells.Item(89, 5).Value = 4.646944933842074 / 100000000
$ws.Cells.Item(89, 6).Value = 3
$ws.Cells.Item(89, 7).Value = "Teta (Duke)"
$ws.Cells.Item(89, 8).Value = 5

$ws.Cells.Item(90, 5).Value = 0.05010556668777629

# rows 92 & 137: "Shaak Ti" (row 92) and "Sheltay Retrac" (row 137) swap places
$ws.Cells.Item(92, 2).Value = 0
$ws.Cells.Item(92, 3).Value = 1
$ws.Cells.Item(92, 4).Value = 0.007042253521126761
$ws.Cells.Item(92, 5).Value = 1.77266881326373 / 10000000000000
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = "Sheltay Retrac"
$ws.Cells.Item(92, 8).Value = 1

$ws.Cells.Item(96, 5).Value = 0.005073025569128987
$ws.Cells.Item(97, 5).Value = 0.1247018671966885
$ws.Cells.Item(98, 5).Value = 0.0007696493910233482
$ws.Cells.Item(102, 5).Value = 0.08530415169585871
$ws.Cells.Item(110, 5).Value = 0.04608679964549955
$ws.Cells.Item(114, 5).Value = 5.233627404279839 / 1000000000
$ws.Cells.Item(124, 5).Value = 0.00711036195663156
$ws.Cells.Item(135, 2).Value = 0.002243414834164804
$ws.Cells.Item(135, 5).Value = 0.1686111444359755
$ws.Cells.Item(136, 5).Value = 0.09083168374410448

$ws.Cells.Item(137, 2).Value = 0.0007169374581189485
$ws.Cells.Item(137, 3).Value = 29
$ws.Cells.Item(137, 4).Value = 0.2042253521126761
$ws.Cells.Item(137, 5).Value = 0.103144667348067
$ws.Cells.Item(137, 6).Value = 12
$ws.Cells.Item(137, 7).Value = "Shaak Ti"
$ws.Cells.Item(137, 8).Value = 17

$ws.Cells.Item(141, 5).Value = 0.02449557211255044
